$wb = $excel.ActiveWorkbook

# Overview sheet: refreshed "Latest HO Xliff Generate Date" for the
# 289b639f-... row (row 2) as part of (re-)generating the handback report.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-12 20:59:04"

# zh-cn sheet: refreshed Correspond Handoff / Handback datetimes for the
# 289b639f-... row (row 2) as part of generating the handback report.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("H2").Value = "2016-08-12 20:58:54"
$wsZh.Range("K2").Value = "2016-08-12 20:59:25"

# de-de sheet: same refresh for the 289b639f-... row (row 2).
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("H2").Value = "2016-08-12 20:59:04"
$wsDe.Range("K2").Value = "2016-08-12 20:59:35"
